# feat: add 2022-Q1 data
#
# 1. Insert a new worksheet "2022-Q1" right before the "总计" (total) sheet,
#    cloning the layout/styling of the most recent per-fund sheet ("2021-Q4")
#    and filling in the 2022-Q1 fund holdings.
# 2. Prepend a new "2022-Q1" row to the "总计" summary sheet, pushing the
#    existing history rows down by one and renumbering the index column.

$wb = $excel.ActiveWorkbook

$template = $wb.Worksheets.Item("2021-Q4")

# --- 1. New "2022-Q1" sheet, inserted immediately before "总计" ---------------
# NOTE: sheet object references resolve by position, so inserting a sheet
# shifts the index every other worksheet reference was bound to. Grab a
# fresh "总计" handle right before the insert, and re-resolve it again
# afterwards (by name) rather than reusing the pre-insert reference.
$totalSheetBefore = $wb.Worksheets.Item("总计")
$newSheet = $wb.Worksheets.Add($totalSheetBefore)
$newSheet.Name = "2022-Q1"

# Clone header/style/layout (incl. the numeric-typed H column) from 2021-Q4.
$template.Range("A1:H3").Copy($newSheet.Range("A1:H3"))

# Header row (unchanged from the template, restated for clarity).
$newSheet.Range("B1").Value = "基金代码"
$newSheet.Range("C1").Value = "基金名称"
$newSheet.Range("D1").Value = "基金规模"
$newSheet.Range("E1").Value = "股票总仓位"
$newSheet.Range("F1").Value = "仓位占比"
$newSheet.Range("G1").Value = "持有市值(亿元)"
$newSheet.Range("H1").Value = "仓位排名"

# Row 2 - 广发纳斯达克生物科技指数(QDII)（人民币）
$newSheet.Range("A2").Value = 0
$newSheet.Range("B2").Value = "'001092"
$newSheet.Range("B2").Style = "Normal"
$newSheet.Range("C2").Value = "广发纳斯达克生物科技指数(QDII)（人民币）"
$newSheet.Range("D2").Value = "'1.34"
$newSheet.Range("D2").Style = "Normal"
$newSheet.Range("E2").Value = "'82.00"
$newSheet.Range("E2").Style = "Normal"
$newSheet.Range("F2").Value = "'5.63"
$newSheet.Range("F2").Style = "Normal"
$newSheet.Range("G2").Value = "'0.0754"
$newSheet.Range("G2").Style = "Normal"
$newSheet.Range("H2").Value = 3

# Row 3 - 广发纳斯达克生物科技指数(QDII)（美元）
$newSheet.Range("A3").Value = 1
$newSheet.Range("B3").Value = "'001093"
$newSheet.Range("B3").Style = "Normal"
$newSheet.Range("C3").Value = "广发纳斯达克生物科技指数(QDII)（美元）"
$newSheet.Range("D3").Value = "'1.34"
$newSheet.Range("D3").Style = "Normal"
$newSheet.Range("E3").Value = "'82.00"
$newSheet.Range("E3").Style = "Normal"
$newSheet.Range("F3").Value = "'5.63"
$newSheet.Range("F3").Style = "Normal"
$newSheet.Range("G3").Value = "'0.0754"
$newSheet.Range("G3").Style = "Normal"
$newSheet.Range("H3").Value = 3

# --- 2. Push a new row into "总计" summary sheet -----------------------------
# Re-resolve "总计" by name now that the sheet collection has changed shape.
$totalSheet = $wb.Worksheets.Item("总计")

# Copy row 6 down to row 7 first so the last row keeps its A-column styling
# when we extend the table by one row.
$totalSheet.Range("A6:D6").Copy($totalSheet.Range("A7:D7"))

$totalSheet.Range("A2").Value = 0
$totalSheet.Range("B2").Value = "2022-Q1"
$totalSheet.Range("C2").Value = 2
$totalSheet.Range("D2").Value = 0.15

$totalSheet.Range("A3").Value = 1
$totalSheet.Range("B3").Value = "2021-Q4"
$totalSheet.Range("C3").Value = 2
$totalSheet.Range("D3").Value = 0.13

$totalSheet.Range("A4").Value = 2
$totalSheet.Range("B4").Value = "2021-Q3"
$totalSheet.Range("C4").Value = 2
$totalSheet.Range("D4").Value = 0.09

$totalSheet.Range("A5").Value = 3
$totalSheet.Range("B5").Value = "2021-Q2"
$totalSheet.Range("C5").Value = 4
$totalSheet.Range("D5").Value = 0.1

$totalSheet.Range("A6").Value = 4
$totalSheet.Range("B6").Value = "2021-Q1"
$totalSheet.Range("C6").Value = 2
$totalSheet.Range("D6").Value = 0.07000000000000001

$totalSheet.Range("A7").Value = 5
$totalSheet.Range("B7").Value = "2020-Q4"
$totalSheet.Range("C7").Value = 2
$totalSheet.Range("D7").Value = 0.08
